$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 change
$ws.Range("Q3").Value = 1.08

# Row 4 changes
$ws.Range("F4").Value  = 1.01
$ws.Range("G4").Value  = 1.02
$ws.Range("H4").Value  = 150
$ws.Range("I4").Value  = 390
$ws.Range("J4").Value  = 60
$ws.Range("K4").Value  = 900
$ws.Range("L4").Value  = 0
$ws.Range("M4").Value  = 0
$ws.Range("N4").Value  = 0
$ws.Range("O4").Value  = 0
$ws.Range("P4").Value  = 0
$ws.Range("Q4").Value  = 0
$ws.Range("R4").Value  = 2.84
$ws.Range("S4").Value  = 1.53
$ws.Range("T4").Value  = 2.28
$ws.Range("U4").Value  = 1.67
$ws.Range("V4").Value  = 1.01
$ws.Range("W4").Value  = 50
$ws.Range("X4").Value  = 1000
$ws.Range("Y4").Value  = 1000
$ws.Range("Z4").Value  = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 1000
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 1000
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 2.76
$ws.Range("AK4").Value = 5.1
$ws.Range("AL4").Value = 19.5
$ws.Range("AM4").Value = 1000
$ws.Range("AN4").Value = 2.82
$ws.Range("AO4").Value = 1000
